$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped coinranking.com snapshot: Price (D) / Volume(1h) (E)
# cells for nearly every row, plus a Chainlink/BinanceUSD row-order swap
# (rows 23-24) that also changed each row's Coin/Link/Price/Volume.
#
# Column D holds "Price" as plain text (note the "." thousands separators,
# e.g. "25.912.66") in the source file. Several new values (e.g. "4.340",
# "6.050") look like ordinary decimals, and Excel's Range.Value setter would
# silently reinterpret them as numbers and drop the trailing zero (4.340 ->
# 4.34). Force the cell to Text first so the literal string is preserved,
# then reset the style to Normal right after so no stray cell format is left
# behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.912.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5163"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06212"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07542"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.645.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.862.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5397"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.947.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.606"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "183.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.942"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.55%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.050"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.294"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("E27").Value = "  -5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05899"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.328"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.317"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9631"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.383"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5799"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01589"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8385"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.031.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.661"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.788.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9962"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.967"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05185"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4232"
$ws.Range("D51").Style = "Normal"
